# Fill in the missing "ü" (checkmark) marks for a handful of students that
# had been left blank in the "Wanted List" monthly task tracker.
#
# Each target cell currently has an empty inline string with a "blank"
# cell style (26 for inner C/D/E/F columns, 27 for the G column). The
# fix sets the value to the checkmark glyph "ü" and restyles the cell to
# match its sibling "checked" cells (style 28 for C/D/E/F, style 33 for
# the G column), by copying formatting from an already-correct neighbour
# cell that uses the desired style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E5: blank -> checked (style 26 -> 28)
$ws.Range("E5").Value = "ü"
$ws.Range("C6").Copy()
$ws.Range("E5").PasteSpecial(-4122)

# E15: blank -> checked (style 26 -> 28)
$ws.Range("E15").Value = "ü"
$ws.Range("C6").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# E20: blank -> checked (style 26 -> 28)
$ws.Range("E20").Value = "ü"
$ws.Range("C6").Copy()
$ws.Range("E20").PasteSpecial(-4122)

# C35: blank -> checked (style 26 -> 28)
$ws.Range("C35").Value = "ü"
$ws.Range("C6").Copy()
$ws.Range("C35").PasteSpecial(-4122)

# G36: blank -> checked (style 27 -> 33), matches other checked G-column cells
$ws.Range("G36").Value = "ü"
$ws.Range("G26").Copy()
$ws.Range("G36").PasteSpecial(-4122)
